$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId 1) - column F ("想去人数") updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 29
$ws1.Range("F6").Value = 558
$ws1.Range("F7").Value = 1748
$ws1.Range("F10").Value = 139
$ws1.Range("F11").Value = 1912
$ws1.Range("F13").Value = 166
$ws1.Range("F14").Value = 440
$ws1.Range("F15").Value = 9
$ws1.Range("F19").Value = 19
$ws1.Range("F23").Value = 1029
$ws1.Range("F28").Value = 282

# Sheet "全部类型" (sheetId 4) - column F ("想去人数") updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 29
$ws4.Range("F6").Value = 558
$ws4.Range("F7").Value = 1748
$ws4.Range("F11").Value = 139
$ws4.Range("F12").Value = 1912
$ws4.Range("F14").Value = 166
$ws4.Range("F15").Value = 440
$ws4.Range("F16").Value = 9
$ws4.Range("F20").Value = 19
$ws4.Range("F24").Value = 1029
$ws4.Range("F29").Value = 282
